$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 ("New York -- New York"): the re-scrape for this run hit a GitHub
# API rate limit, so every field that depended on the fetch came back
# blank and the status column records the new error instead of "Success!".
# Clear the now-empty fields (content + restore the default style, since
# the original date cell B4 carried a date-format style that no longer
# applies to a blank cell) and flip the "includes Hispanic Black" boolean
# to False.
$ws.Range("B4:H4").ClearContents()
$ws.Range("B4:H4").Style = "Normal"

$ws.Range("J4").Value = $false

$ws.Range("K4:L4").ClearContents()
$ws.Range("K4:L4").Style = "Normal"

$ws.Range("O4").Value = "An error occurred. ... RateLimitExceededException(403, {'message': ""API rate limit exceeded for 132.145.200.60. (But here's the good news: Authenticated requests get a higher rate limit. Check out the documentation for more details.)"", 'documentation_url': 'https://developer.github.com/v3/#rate-limiting'})"

# Row 39 (Delaware): same run, different failure mode for this source --
# previously a 504 gateway timeout, now a numpy AttributeError further
# along in the parsing step.
$ws.Range("O39").Value = "An error occurred. ... AttributeError(""'numpy.float64' object has no attribute 'split'"")"
